# Add a new diary entry for Aman Bhatia (diary-AmanBhatia.xlsx)
#
# This adds three new diary rows (13, 14, 15) covering 2020-01-16,
# 2020-01-18 and 2020-01-20, formatted like the existing entries in
# rows 10-12, and updates the Participants cell of the "Team formation"
# entry (row 11) to include "Myself".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare rows 13-15 with the same look & feel as rows 10-12 -----------
# Copy the formatting (number format / font / fill / wrap / etc.) of row 10
# (A:G) down onto the three new rows before filling in their content.
$ws.Range("A10:G10").Copy()
$ws.Range("A13:G13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A10:G10").Copy()
$ws.Range("A14:G14").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A10:G10").Copy()
$ws.Range("A15:G15").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 13 : 2020-01-16 ----------------------------------------------------
$ws.Cells.Item(13, 1).Value = 43846
$ws.Cells.Item(13, 2).Value = "5:00 - 7:50 pm"
$ws.Cells.Item(13, 3).Value = "Myself"
$ws.Cells.Item(13, 4).Value = "To understand the technicalities of Reverse Engineering, by actually analyzing code, and try to find the causes of bugs."
$ws.Cells.Item(13, 6).Value = 'Reverse engineering is not easy, but a little understanding can help us in finding beacons, which can guide us to the correct answer, in todays lecture, this correct answer was the location of the bugs/feature'
$ws.Cells.Item(13, 5).Value = "We were able to find more or less the precise causes of bugs, in JPacMan 1 and 2, by finding usages of relevant classes,methods,keywords.Used various approaches like opportunist,bottom up etc."
$ws.Cells.Item(13, 7).Value = "excited"
$ws.Rows.Item(13).RowHeight = 102

# --- Row 14 : 2020-01-18 ----------------------------------------------------
$ws.Cells.Item(14, 1).Value = 43848
$ws.Cells.Item(14, 2).Value = "6:00 - 10:00pm"
$ws.Cells.Item(14, 3).Value = "Anjana, Vaishakhi,Myself"
$ws.Cells.Item(14, 4).Value = "To decide on a group projct"
$ws.Cells.Item(14, 5).Value = "After reading a lot of projects, we finally found a relevant one -OpenRefine"
$ws.Cells.Item(14, 6).Value = 'It is very difficult to find the "perfect" project, some of them have cryptic looking code, but a rich documentation, while others lack a documentation'
$ws.Cells.Item(14, 7).Value = "Neutral"
$ws.Rows.Item(14).RowHeight = 68

# --- Row 15 : 2020-01-20 ----------------------------------------------------
$ws.Cells.Item(15, 1).Value = 43850
$ws.Cells.Item(15, 2).Value = "10:00 am - 12:30 pm"
$ws.Cells.Item(15, 3).Value = "Myself"
$ws.Cells.Item(15, 4).Value = "To finsih the homework, and understand how we could add new features to an existing codebase"
$ws.Cells.Item(15, 5).Value = "Was able to understand what the homework asked of me. Explored JPacMan3, and understood, the process of building up on already existing code."
$ws.Cells.Item(15, 6).Value = "If the code employs good coding practices, like proper nomenclature of variables, proper method names, comments, it. Becomes slightly easier to understand code, even if somebody else had written it. The third question was interesting, since I felt, it could be done in more than one way"
$ws.Cells.Item(15, 7).Value = "Neutral"
$ws.Rows.Item(15).RowHeight = 136

# --- Update existing row 11 Participants cell (reuses the string created
#     above for C14) -----------------------------------------------------
$ws.Cells.Item(11, 3).Value = "Anjana, Vaishakhi,Myself"

# --- Update the sheet view: scrolled to row 10, selection on G15 ----------
$ws.Range("G15").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
